$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2654  # Afghanistan: 2642 -> 2654
$ws.Range("B3").Value = 2403  # Albania: 2399 -> 2403
$ws.Range("B4").Value = 3299  # Algeria: 3280 -> 3299
$ws.Range("B6").Value = 618  # Angola: 603 -> 618
$ws.Range("B8").Value = 65865  # Argentina: 64792 -> 65865
$ws.Range("B9").Value = 4178  # Armenia: 4149 -> 4178
$ws.Range("B11").Value = 10311  # Austria: 10260 -> 10311
$ws.Range("B12").Value = 4617  # Azerbaijan: 4580 -> 4617
$ws.Range("B14").Value = 664  # Bahrain: 657 -> 664
$ws.Range("B15").Value = 11755  # Bangladesh: 11644 -> 11755
$ws.Range("B16").Value = 45  # Barbados: 44 -> 45
$ws.Range("B17").Value = 2592  # Belarus: 2572 -> 2592
$ws.Range("B18").Value = 24406  # Belgium: 24322 -> 24406
$ws.Range("B22").Value = 13082  # Bolivia: 13035 -> 13082
$ws.Range("B23").Value = 8736  # Bosnia and Herzegovina: 8551 -> 8736
$ws.Range("B25").Value = 414399  # Brazil: 408622 -> 414399
$ws.Range("B27").Value = 16773  # Bulgaria: 16548 -> 16773
$ws.Range("B28").Value = 161  # Burkina Faso: 157 -> 161
$ws.Range("B29").Value = 3210  # Burma: 3209 -> 3210
$ws.Range("B31").Value = 229  # Cabo Verde: 226 -> 229
$ws.Range("B32").Value = 110  # Cambodia: 106 -> 110
$ws.Range("B34").Value = 24437  # Canada: 24334 -> 24437
$ws.Range("B35").Value = 91  # Central African Republic: 88 -> 91
$ws.Range("B37").Value = 26726  # Chile: 26659 -> 26726
$ws.Range("B39").Value = 76015  # Colombia: 75164 -> 76015
$ws.Range("B42").Value = 770  # Congo (Kinshasa): 768 -> 770
$ws.Range("B43").Value = 3326  # Costa Rica: 3290 -> 3326
$ws.Range("B44").Value = 291  # Cote d'Ivoire: 287 -> 291
$ws.Range("B45").Value = 7315  # Croatia: 7218 -> 7315
$ws.Range("B46").Value = 694  # Cuba: 675 -> 694
$ws.Range("B47").Value = 327  # Cyprus: 321 -> 327
$ws.Range("B48").Value = 29479  # Czechia: 29365 -> 29479
$ws.Range("B49").Value = 2492  # Denmark: 2491 -> 2492
$ws.Range("B50").Value = 147  # Djibouti: 146 -> 147
$ws.Range("B52").Value = 3509  # Dominican Republic: 3499 -> 3509
$ws.Range("B53").Value = 18907  # Ecuador: 18765 -> 18907
$ws.Range("B54").Value = 13655  # Egypt: 13531 -> 13655
$ws.Range("B55").Value = 2141  # El Salvador: 2134 -> 2141
$ws.Range("B58").Value = 1183  # Estonia: 1172 -> 1183
$ws.Range("B60").Value = 3795  # Ethiopia: 3757 -> 3795
$ws.Range("B62").Value = 918  # Finland: 915 -> 918
$ws.Range("B63").Value = 105792  # France: 105291 -> 105792
$ws.Range("B66").Value = 4207  # Georgia: 4163 -> 4207
$ws.Range("B67").Value = 84141  # Germany: 83605 -> 84141
$ws.Range("B68").Value = 783  # Ghana: 779 -> 783
$ws.Range("B69").Value = 10764  # Greece: 10587 -> 10764
$ws.Range("B71").Value = 7642  # Guatemala: 7578 -> 7642
$ws.Range("B72").Value = 148  # Guinea: 146 -> 148
$ws.Range("B74").Value = 308  # Guyana: 303 -> 308
$ws.Range("B75").Value = 260  # Haiti: 254 -> 260
$ws.Range("B76").Value = 5439  # Honduras: 5367 -> 5439
$ws.Range("B77").Value = 28173  # Hungary: 27908 -> 28173
$ws.Range("B79").Value = 230168  # India: 222408 -> 230168
$ws.Range("B80").Value = 46349  # Indonesia: 45949 -> 46349
$ws.Range("B81").Value = 73568  # Iran: 72875 -> 73568
$ws.Range("B82").Value = 15640  # Iraq: 15566 -> 15640
$ws.Range("B83").Value = 4915  # Ireland: 4906 -> 4915
$ws.Range("B84").Value = 6370  # Israel: 6367 -> 6370
$ws.Range("B85").Value = 122005  # Italy: 121433 -> 122005
$ws.Range("B86").Value = 794  # Jamaica: 790 -> 794
$ws.Range("B87").Value = 10501  # Japan: 10391 -> 10501
$ws.Range("B88").Value = 8988  # Jordan: 8925 -> 8988
$ws.Range("B89").Value = 3363  # Kazakhstan: 3352 -> 3363
$ws.Range("B90").Value = 2825  # Kenya: 2781 -> 2825
$ws.Range("B91").Value = 1851  # Korea, South: 1840 -> 1851
$ws.Range("B93").Value = 1610  # Kuwait: 1590 -> 1610
$ws.Range("B94").Value = 1637  # Kyrgyzstan: 1622 -> 1637
$ws.Range("B96").Value = 2166  # Latvia: 2154 -> 2166
$ws.Range("B97").Value = 7390  # Lebanon: 7345 -> 7390
$ws.Range("B98").Value = 319  # Lesotho: 318 -> 319
$ws.Range("B100").Value = 3058  # Libya: 3047 -> 3058
$ws.Range("B102").Value = 3993  # Lithuania: 3956 -> 3993
$ws.Range("B104").Value = 690  # Madagascar: 677 -> 690
$ws.Range("B105").Value = 1151  # Malawi: 1148 -> 1151
$ws.Range("B106").Value = 1591  # Malaysia: 1551 -> 1591
$ws.Range("B108").Value = 493  # Mali: 491 -> 493
$ws.Range("B111").Value = 456  # Mauritania: 455 -> 456
$ws.Range("B113").Value = 218007  # Mexico: 217345 -> 218007
$ws.Range("B114").Value = 5892  # Moldova: 5850 -> 5892
$ws.Range("B116").Value = 130  # Mongolia: 126 -> 130
$ws.Range("B117").Value = 1521  # Montenegro: 1510 -> 1521
$ws.Range("B118").Value = 9043  # Morocco: 9032 -> 9043
$ws.Range("B119").Value = 818  # Mozambique: 815 -> 818
$ws.Range("B120").Value = 667  # Namibia: 643 -> 667
$ws.Range("B121").Value = 3475  # Nepal: 3362 -> 3475
$ws.Range("B122").Value = 17501  # Netherlands: 17443 -> 17501
$ws.Range("B124").Value = 183  # Nicaragua: 182 -> 183
$ws.Range("B125").Value = 192  # Niger: 191 -> 192
$ws.Range("B126").Value = 2065  # Nigeria: 2063 -> 2065
$ws.Range("B127").Value = 767  # Norway: 757 -> 767
$ws.Range("B128").Value = 2071  # Oman: 2053 -> 2071
$ws.Range("B129").Value = 18429  # Pakistan: 18310 -> 18429
$ws.Range("B130").Value = 6252  # Panama: 6244 -> 6252
$ws.Range("B131").Value = 121  # Papua New Guinea: 115 -> 121
$ws.Range("B132").Value = 6798  # Paraguay: 6653 -> 6798
$ws.Range("B133").Value = 62674  # Peru: 62375 -> 62674
$ws.Range("B134").Value = 17800  # Philippines: 17525 -> 17800
$ws.Range("B135").Value = 68482  # Poland: 68105 -> 68482
$ws.Range("B136").Value = 16983  # Portugal: 16977 -> 16983
$ws.Range("B137").Value = 489  # Qatar: 480 -> 489
$ws.Range("B138").Value = 28616  # Romania: 28380 -> 28616
$ws.Range("B139").Value = 110022  # Russia: 109341 -> 110022
$ws.Range("B140").Value = 338  # Rwanda: 337 -> 338
$ws.Range("B147").Value = 7018  # Saudi Arabia: 6992 -> 7018
$ws.Range("B148").Value = 1114  # Senegal: 1111 -> 1114
$ws.Range("B149").Value = 6478  # Serbia: 6432 -> 6478
$ws.Range("B153").Value = 11886  # Slovakia: 11807 -> 11886
$ws.Range("B154").Value = 4279  # Slovenia: 4269 -> 4279
$ws.Range("B156").Value = 721  # Somalia: 713 -> 721
$ws.Range("B157").Value = 54557  # South Africa: 54452 -> 54557
$ws.Range("B159").Value = 78566  # Spain: 78293 -> 78566
$ws.Range("B160").Value = 734  # Sri Lanka: 709 -> 734
$ws.Range("B161").Value = 2365  # Sudan: 2349 -> 2365
$ws.Range("B162").Value = 212  # Suriname: 207 -> 212
$ws.Range("B163").Value = 14151  # Sweden: 14048 -> 14151
$ws.Range("B164").Value = 10685  # Switzerland: 10655 -> 10685
$ws.Range("B165").Value = 1625  # Syria: 1610 -> 1625
$ws.Range("B169").Value = 318  # Thailand: 276 -> 318
$ws.Range("B172").Value = 189  # Trinidad and Tobago: 179 -> 189
$ws.Range("B173").Value = 11122  # Tunisia: 10915 -> 11122
$ws.Range("B174").Value = 41883  # Turkey: 41191 -> 41883
$ws.Range("B175").Value = 579275  # US: 577528 -> 579275
$ws.Range("B176").Value = 343  # Uganda: 342 -> 343
$ws.Range("B177").Value = 46950  # Ukraine: 46607 -> 46950
$ws.Range("B178").Value = 1601  # United Arab Emirates: 1596 -> 1601
$ws.Range("B179").Value = 127830  # United Kingdom: 127797 -> 127830
$ws.Range("B180").Value = 2918  # Uruguay: 2796 -> 2918
$ws.Range("B181").Value = 655  # Uzbekistan: 653 -> 655
$ws.Range("B183").Value = 2226  # Venezuela: 2189 -> 2226
$ws.Range("B185").Value = 3317  # West Bank and Gaza: 3283 -> 3317
$ws.Range("B186").Value = 1260  # Yemen: 1239 -> 1260
$ws.Range("B187").Value = 1255  # Zambia: 1253 -> 1255
$ws.Range("B188").Value = 1574  # Zimbabwe: 1573 -> 1574
